$p = $ppt.ActivePresentation

# Locate the "Redundancies" slide (slide 6 in the deck) by its title so the
# script is resilient even if slide ordering assumptions are ever wrong.
$slide = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    if ($candidate.Shapes.Item(1).TextFrame.TextRange.Text -eq "Redundancies") {
        $slide = $candidate
        break
    }
}
if ($slide -eq $null) {
    $slide = $p.Slides.Item(6)
}

$shape = $slide.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# --- 1. Insert a new sub-bullet "Downside ... entrypoint.sh" before the
#        "The use of storage containers" paragraph (currently paragraph 2). ---
$para2 = $tr.Paragraphs(2)
$para2.InsertBefore("Downside – potentially huge entrypoint.sh `r")

# Re-fetch the range/paragraphs after the structural edit and demote the
# freshly inserted paragraph to outline level 2 (lvl="1" in the XML).
$tr = $shape.TextFrame.TextRange
$tr.Paragraphs(2).IndentLevel = 2

# --- 2. Update the (now 3rd) paragraph text to add the "(No longer
#        recommended)" suffix. Assign an unrelated placeholder first so the
#        host doesn't keep the old text as a separate shared-prefix run -
#        this keeps the final paragraph to a single run. ---
$tr = $shape.TextFrame.TextRange
$para3 = $tr.Paragraphs(3)
$para3.Text = "~tmp~"
$tr = $shape.TextFrame.TextRange
$tr.Paragraphs(3).Text = "The use of storage containers (No longer recommended)"

# --- 3. Merge the two runs of the final paragraph ("No data is
#        theoretically lost should a container " + "go down.") into a
#        single run with the same visible text. ---
$tr = $shape.TextFrame.TextRange
$para5 = $tr.Paragraphs(5)
$para5.Text = "~tmp~"
$tr = $shape.TextFrame.TextRange
$tr.Paragraphs(5).Text = "No data is theoretically lost should a container go down."
